$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I2").Value = 0.9223046214701632
$ws.Range("J2").Value = 0.9223046214701632
$ws.Range("M2").Value = 4.093680666666667
$ws.Range("N2").Value = 12.281042
$ws.Range("O2").Value = 0.1610908176055751
$ws.Range("P2").Value = 0.161090817605575
$ws.Range("Q2").Value = 6.830392159627334
$ws.Range("R2").Value = 61.473529436646
$ws.Range("S2").Value = 0.148574805554029
$ws.Range("T2").Value = 0.148574805554029
$ws.Range("I3").Value = 0.9223046214701632
$ws.Range("J3").Value = 0.9223046214701632
$ws.Range("O3").Value = 0.5606512265211691
$ws.Range("P3").Value = 0.5606512265211691
$ws.Range("S3").Value = 0.5170912172533896
$ws.Range("T3").Value = 0.5170912172533896
$ws.Range("I4").Value = 0.9223046214701632
$ws.Range("J4").Value = 0.9223046214701632
$ws.Range("M4").Value = 7.071161666666666
$ws.Range("N4").Value = 21.213485
$ws.Range("O4").Value = 0.2782579558732559
$ws.Range("P4").Value = 0.2782579558732559
$ws.Range("Q4").Value = 11.79838173522833
$ws.Range("R4").Value = 106.185435617055
$ws.Range("S4").Value = 0.2566385986627446
$ws.Range("T4").Value = 0.2566385986627446
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0.140557
$ws.Range("H5").Value = 0.421671
$ws.Range("I5").Value = 0.07769537852983674
$ws.Range("J5").Value = 0.07769537852983674
$ws.Range("M5").Value = 4.093680666666667
$ws.Range("N5").Value = 12.281042
$ws.Range("O5").Value = 0.1610908176055751
$ws.Range("P5").Value = 0.161090817605575
$ws.Range("Q5").Value = 0.5753954734646668
$ws.Range("R5").Value = 5.178559261182
$ws.Range("S5").Value = 0.01251601205154604
$ws.Range("T5").Value = 0.01251601205154604
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 0.140557
$ws.Range("H6").Value = 0.421671
$ws.Range("I6").Value = 0.07769537852983674
$ws.Range("J6").Value = 0.07769537852983674
$ws.Range("O6").Value = 0.5606512265211691
$ws.Range("P6").Value = 0.5606512265211691
$ws.Range("Q6").Value = 2.002573347927
$ws.Range("R6").Value = 18.023160131343
$ws.Range("S6").Value = 0.04356000926777948
$ws.Range("T6").Value = 0.04356000926777948
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 0.140557
$ws.Range("H7").Value = 0.421671
$ws.Range("I7").Value = 0.07769537852983674
$ws.Range("J7").Value = 0.07769537852983674
$ws.Range("M7").Value = 7.071161666666666
$ws.Range("N7").Value = 21.213485
$ws.Range("O7").Value = 0.2782579558732559
$ws.Range("P7").Value = 0.2782579558732559
$ws.Range("Q7").Value = 0.9939012703816666
$ws.Range("R7").Value = 8.945111433435001
$ws.Range("S7").Value = 0.02161935721051122
$ws.Range("T7").Value = 0.02161935721051122